# Auto-generated Excel COM-interop edit script
# Applies the numeric cell updates described in the commit diff
# across sheets ALC, ARM, BSM, CUL, WVR.

$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H45").Value = 3000
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 3000
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 9000
$ws.Range("M45").ClearContents()
$ws.Range("N45").Value = -9384
$ws.Range("H63").Value = 38000
$ws.Range("J63").Value = 38000
$ws.Range("L63").Value = 38000
$ws.Range("N63").Value = -39248
$ws.Range("H64").Value = 3735.8696
$ws.Range("I64").Value = 3147.9614
$ws.Range("J64").Value = 5534.1763
$ws.Range("K64").Value = 3147.9614
$ws.Range("L64").Value = 5534.1763
$ws.Range("M64").Value = -2899.9614
$ws.Range("N64").Value = -6030.1763
$ws.Range("H66").Value = 38000
$ws.Range("J66").Value = 38000
$ws.Range("L66").Value = 114000
$ws.Range("N66").Value = -120240
$ws.Range("H67").Value = 3735.8696
$ws.Range("I67").Value = 3147.9614
$ws.Range("J67").Value = 5534.1763
$ws.Range("K67").Value = 3147.9614
$ws.Range("L67").Value = 5534.1763
$ws.Range("M67").Value = -2289.9614
$ws.Range("N67").Value = -7250.1763
$ws.Range("H68").Value = 42000
$ws.Range("J68").Value = 42000
$ws.Range("L68").Value = 42000
$ws.Range("N68").Value = -43498
$ws.Range("H69").Value = 3805.818
$ws.Range("J69").Value = 4354.6924
$ws.Range("L69").Value = 13064.0772
$ws.Range("N69").Value = -14812.0772
$ws.Range("H70").Value = 5300
$ws.Range("I70").Value = 1275
$ws.Range("J70").Value = 9900
$ws.Range("K70").Value = 3825
$ws.Range("L70").Value = 29700
$ws.Range("M70").Value = -3555
$ws.Range("N70").Value = -30240
$ws.Range("H71").Value = 42000
$ws.Range("J71").Value = 42000
$ws.Range("L71").Value = 126000
$ws.Range("N71").Value = -133488
$ws.Range("H72").Value = 3805.818
$ws.Range("J72").Value = 4354.6924
$ws.Range("L72").Value = 39192.2316
$ws.Range("N72").Value = -47928.2316
$ws.Range("H73").Value = 5300
$ws.Range("I73").Value = 1275
$ws.Range("J73").Value = 9900
$ws.Range("K73").Value = 3825
$ws.Range("L73").Value = 29700
$ws.Range("M73").Value = -2889
$ws.Range("N73").Value = -31572
$ws.Range("H137").Value = 2112.318
$ws.Range("I137").Value = 1470.4445
$ws.Range("K137").Value = 4411.333500000001
$ws.Range("M137").Value = -1861.333500000001

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3274.92
$ws.Range("I132").Value = 1695.5834
$ws.Range("K132").Value = 5086.7502
$ws.Range("M132").Value = -2556.7502

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H114").Value = 32000
$ws.Range("J114").Value = 32000
$ws.Range("L114").Value = 32000
$ws.Range("N114").Value = -40678

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 3999
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 3999
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 11997
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -13369
$ws.Range("H63").Value = 5475.25
$ws.Range("I63").Value = 1912
$ws.Range("K63").Value = 5736
$ws.Range("M63").Value = -4987
$ws.Range("H64").Value = 3000
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()
$ws.Range("H65").Value = 3999
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 3999
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 35991
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -42855
$ws.Range("H66").Value = 5475.25
$ws.Range("I66").Value = 1912
$ws.Range("K66").Value = 17208
$ws.Range("M66").Value = -13464
$ws.Range("H67").Value = 3000
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()
$ws.Range("H68").Value = 8000500
$ws.Range("I68").Value = 16000000
$ws.Range("J68").Value = 1000
$ws.Range("K68").Value = 48000000
$ws.Range("L68").Value = 3000
$ws.Range("M68").Value = -47999189
$ws.Range("N68").Value = -4622
$ws.Range("H69").Value = 492.77777
$ws.Range("I69").Value = 520
$ws.Range("J69").Value = 479.16666
$ws.Range("K69").Value = 1560
$ws.Range("L69").Value = 1437.49998
$ws.Range("M69").Value = -749
$ws.Range("N69").Value = -3059.49998
$ws.Range("H70").Value = 2287
$ws.Range("I70").Value = 930.5
$ws.Range("K70").Value = 2791.5
$ws.Range("M70").Value = -2476.5
$ws.Range("H71").Value = 8000500
$ws.Range("I71").Value = 16000000
$ws.Range("J71").Value = 1000
$ws.Range("K71").Value = 144000000
$ws.Range("L71").Value = 9000
$ws.Range("M71").Value = -143995944
$ws.Range("N71").Value = -17112
$ws.Range("H72").Value = 492.77777
$ws.Range("I72").Value = 520
$ws.Range("J72").Value = 479.16666
$ws.Range("K72").Value = 4680
$ws.Range("L72").Value = 4312.49994
$ws.Range("M72").Value = -624
$ws.Range("N72").Value = -12424.49994
$ws.Range("H73").Value = 2287
$ws.Range("I73").Value = 930.5
$ws.Range("K73").Value = 2791.5
$ws.Range("M73").Value = -1699.5
$ws.Range("H74").Value = 8999
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 8999
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 26997
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -29119
$ws.Range("H77").Value = 8999
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 8999
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 80991
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -91599
$ws.Range("H113").Value = 1250.2916
$ws.Range("J113").Value = 1522.3334
$ws.Range("L113").Value = 4567.0002
$ws.Range("N113").Value = -8907.0002

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 23998
$ws.Range("J80").Value = 23998
$ws.Range("L80").Value = 23998
$ws.Range("N80").Value = -25994
$ws.Range("H83").Value = 23998
$ws.Range("J83").Value = 23998
$ws.Range("L83").Value = 71994
$ws.Range("N83").Value = -81978
$ws.Range("H132").Value = 2410.3103
$ws.Range("I132").Value = 1825.8422
$ws.Range("J132").Value = 3520.8
$ws.Range("K132").Value = 5477.5266
$ws.Range("L132").Value = 10562.4
$ws.Range("M132").Value = -2947.5266
$ws.Range("N132").Value = -15622.4
